$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update email addresses for the three collaborators to the new company emails
$ws.Range("E2").Value = "abuitrago@suraelec.com"
$ws.Range("E3").Value = "wromero@suraelec.com"
$ws.Range("E4").Value = "amanzi@suraelec.com"

# 2. Refresh hyperlinks so the mailto targets match the new addresses
$ws.Range("E2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:abuitrago@suraelec.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:wromero@suraelec.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:amanzi@suraelec.com")
$ws.Range("E2:E4").Style = "Hyperlink"

# 3. Add the new "profile" column to the table
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Add() | Out-Null
$ws.Range("F1").Value = "profile"

# 4. Remove an empty row below the table data (shifts the trailing formatted row up)
$ws.Rows.Item(5).Delete()

# 5. Update active selection
$ws.Range("E16").Select()
